# Commit: Storage: simple cyclic SoC - Formulation adapted from PyPSA -
#         Lack of high VRE nullifies storage - Constant import price nullifies storage
#
# Concretely, for this workbook (CHE_convchp_fueloil), the edit inserts a new
# data row right after the header block (before the existing "input"/"output"
# parameter rows), adding a new "enable_year" / "configuration" parameter
# entry with Value 1990. All subsequent rows shift down by one, and the
# sheet's AutoFilter range / _FilterDatabase defined name grow accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7 (pushes existing row 7.. down to row 8..)
$ws.Rows("7:7").Insert()

# Populate the new row: Country, Entity, Parameter, Type, (Year blank,
# Flow blank), Value, (Unit/Delete/Reference/Link/Note blank)
$ws.Range("A7").Value = "CHE"
$ws.Range("B7").Value = "conv_chp_oil"
$ws.Range("C7").Value = "enable_year"
$ws.Range("D7").Value = "configuration"
$ws.Range("G7").Value = 1990

# Refresh the AutoFilter so its range covers the extra row (5..853)
$ws.AutoFilterMode = $false
$ws.Range("A5:L853").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range
foreach ($nm in $ws.Names) {
    if ($nm.Name -eq "Sheet1!_FilterDatabase") {
        $nm.RefersTo = "=Sheet1!`$A`$5:`$L`$853"
    }
}

# Match the author's last active selection after the edit
$ws.Range("H7").Select()
